$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# Find and delete the shape named "Rectangle 17" located at (1763673, 3442275)
for ($i = $s.Shapes.Count; $i -ge 1; $i--) {
    $shp = $s.Shapes.Item($i)
    if ($shp.Name -eq "Rectangle 17") {
        $shp.Delete()
    }
}
